$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base")
$cfg = $wb.Worksheets.Item("Config")

# --- Base sheet: insert new "CUENTA" column before F_INICIO (old column J) ---
$ws.Columns.Item(10).Insert()

# New header cell, copy formatting from the ESTADO/CLIENTE-style header (B1) which matches the target style
$ws.Range("J1").Value = "CUENTA"
$ws.Range("B1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New CUENTA values for existing rows
$ws.Range("J2").Value = "BULL MARKETING"
$ws.Range("J3").Value = "BULL MARKETING"
$ws.Range("J4").Value = "V2V"

# --- Config sheet: add a second list column (C) with CUENTA options ---
$cfg.Range("C1").Value = "BULL MARKETING"
$cfg.Range("C2").Value = "V2V"
$cfg.Columns.Item(3).ColumnWidth = 13

# --- Base sheet: add the data validation dropdown for the new CUENTA column ---
$validationRange = $ws.Range("J2:J1048576")
$validationRange.Validation.Add(3, 1, 1, "=Config!`$C`$1:`$C`$2")
$validationRange.Validation.InCellDropdown = $true

$wb.Save()
